$d = $word.ActiveDocument

function Merge-ParagraphRuns($doc, $paraIndex) {
    $p = $doc.Paragraphs($paraIndex)
    $full = $p.Range
    $full.End = $full.End - 1
    $paraStart = $full.Start
    $paraEnd = $full.End
    $text = $full.Text

    $spaceIdx = $text.IndexOf(" ")
    if ($spaceIdx -lt 0) {
        return
    }

    $run1End = $paraStart + $spaceIdx
    $tailText = $text.Substring($spaceIdx)

    # Remove everything after the paragraph's first run (first word), then
    # re-insert the remainder right onto the end of that first run so the
    # existing run absorbs the rest of the paragraph's text instead of a
    # brand new run/text node being minted for it.
    $tail = $doc.Range($run1End, $paraEnd)
    $tail.Delete()

    $insertPoint = $doc.Range($run1End, $run1End)
    $insertPoint.InsertAfter($tailText)
}

# Title: "Desmos now in STARMAST"
Merge-ParagraphRuns $d 1

# Author: "Tom Coleman"
Merge-ParagraphRuns $d 2

# Abstract: "Desmos figures now included in STARMAST resources!"
Merge-ParagraphRuns $d 5
